# Appends the new "B2-B2" schedule rows (122-161) to the single
# worksheet, continuing the existing table (Year / Group / Subject / Session /
# Date / Start Time / Duration) and mirroring its alternating row styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Year, Group, Subject, Session, Date, Start Time, Duration(min)
$data = @(
  @('Year 5','B2-B2','endocrinology','1','24/01/2026','09:00:00',360),
  @('Year 5','B2-B2','endocrinology','2','25/01/2026','09:00:00',360),
  @('Year 5','B2-B2','endocrinology','3','26/01/2026','09:00:00',360),
  @('Year 5','B2-B2','endocrinology','4','27/01/2026','09:00:00',360),
  @('Year 5','B2-B2','endocrinology','5','28/01/2026','09:00:00',360),
  @('Year 5','B2-B2','endocrinology','6','07/02/2026','09:00:00',360),
  @('Year 5','B2-B2','endocrinology','7','08/02/2026','09:00:00',360),
  @('Year 5','B2-B2','endocrinology','8','09/02/2026','09:00:00',360),
  @('Year 5','B2-B2','endocrinology','9','10/02/2026','09:00:00',360),
  @('Year 5','B2-B2','endocrinology','10','11/02/2026','09:00:00',360),
  @('Year 5','B2-B2','gastroenterology','1','06/12/2025','09:00:00',360),
  @('Year 5','B2-B2','gastroenterology','2','07/12/2025','09:00:00',360),
  @('Year 5','B2-B2','gastroenterology','3','08/12/2025','09:00:00',360),
  @('Year 5','B2-B2','gastroenterology','4','09/12/2025','09:00:00',360),
  @('Year 5','B2-B2','gastroenterology','5','10/12/2025','09:00:00',360),
  @('Year 5','B2-B2','gastroenterology','6','13/12/2025','09:00:00',360),
  @('Year 5','B2-B2','gastroenterology','7','14/12/2025','09:00:00',360),
  @('Year 5','B2-B2','gastroenterology','8','15/12/2025','09:00:00',360),
  @('Year 5','B2-B2','gastroenterology','9','16/12/2025','09:00:00',360),
  @('Year 5','B2-B2','gastroenterology','10','17/12/2025','09:00:00',360),
  @('Year 5','B2-B2','nephrology','1','20/12/2025','09:00:00',360),
  @('Year 5','B2-B2','nephrology','2','21/12/2025','09:00:00',360),
  @('Year 5','B2-B2','nephrology','3','22/12/2025','09:00:00',360),
  @('Year 5','B2-B2','nephrology','4','23/12/2025','09:00:00',360),
  @('Year 5','B2-B2','nephrology','5','24/12/2025','09:00:00',360),
  @('Year 5','B2-B2','neurology','1','03/01/2026','09:00:00',360),
  @('Year 5','B2-B2','neurology','2','04/01/2026','09:00:00',360),
  @('Year 5','B2-B2','neurology','3','05/01/2026','09:00:00',360),
  @('Year 5','B2-B2','neurology','4','06/01/2026','09:00:00',360),
  @('Year 5','B2-B2','neurology','5','10/01/2026','09:00:00',360),
  @('Year 5','B2-B2','neurology','6','11/01/2026','09:00:00',360),
  @('Year 5','B2-B2','neurology','7','12/01/2026','09:00:00',360),
  @('Year 5','B2-B2','neurology','8','13/01/2026','09:00:00',360),
  @('Year 5','B2-B2','physical medicine','1','07/01/2026','09:00:00',360),
  @('Year 5','B2-B2','physical medicine','2','14/01/2026','09:00:00',360),
  @('Year 5','B2-B2','rheumatology','1','27/12/2025','09:00:00',360),
  @('Year 5','B2-B2','rheumatology','2','28/12/2025','09:00:00',360),
  @('Year 5','B2-B2','rheumatology','3','29/12/2025','09:00:00',360),
  @('Year 5','B2-B2','rheumatology','4','30/12/2025','09:00:00',360),
  @('Year 5','B2-B2','rheumatology','5','31/12/2025','09:00:00',360)
)

$startRow = 122
$endRow   = 161
$lastExistingEvenRow = 120   # reference row using style set s=2,2,2,2,3,4,5 (A..G)
$lastExistingOddRow  = 121   # reference row using style set s=6,6,6,6,7,8,9 (A..G)

# Pre-format the target text columns (A:F) as Text so that numeric-looking
# strings (Session numbers, dates, times) are written as text values rather
# than being auto-converted to numbers/dates - matching the source data,
# which stores every column except Duration as text.
$ws.Range("A" + $startRow + ":F" + $endRow).NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = [string]$row[0]
    $ws.Cells.Item($r, 2).Value = [string]$row[1]
    $ws.Cells.Item($r, 3).Value = [string]$row[2]
    $ws.Cells.Item($r, 4).Value = [string]$row[3]
    $ws.Cells.Item($r, 5).Value = [string]$row[4]
    $ws.Cells.Item($r, 6).Value = [string]$row[5]
    $ws.Cells.Item($r, 7).Value = [double]$row[6]

    # Copy the visual formatting (fill/font/number format/style index) from
    # the matching-parity row already in the table, without touching the
    # values just set (PasteSpecial xlPasteFormats = -4122 copies formats only).
    if (($r % 2) -eq 0) {
        $ws.Range("A" + $lastExistingEvenRow + ":G" + $lastExistingEvenRow).Copy()
    } else {
        $ws.Range("A" + $lastExistingOddRow + ":G" + $lastExistingOddRow).Copy()
    }
    $ws.Range("A" + $r + ":G" + $r).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
